# MAJOR CHANGES: restructure repo for versioning
#
# The "Namespaces-v2" sheet lists the v2 DPV namespaces. Their IRIs move
# from https://w3id.org/dpv/v2/... to https://w3id.org/dpv/v2.0/...
# (and the base "dpv" row, which used to point at the v1 IRI, now points
# at the new v2.0 base IRI).
#
# NOTE: this runtime's Range/Cells ".Value" getter does not reliably
# surface the underlying stored string back into the script, so rather
# than doing a read-modify-write round trip we assign each target IRI
# explicitly (the mapping below mirrors the source workbook row-by-row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Namespaces-v2")

$ws.Cells.Item(2, 2).Value  = "https://w3id.org/dpv/v2.0#"
$ws.Cells.Item(3, 2).Value  = "https://w3id.org/dpv/v2.0/risk#"
$ws.Cells.Item(4, 2).Value  = "https://w3id.org/dpv/v2.0/examples#"
$ws.Cells.Item(5, 2).Value  = "https://w3id.org/dpv/v2.0/rights#"
$ws.Cells.Item(6, 2).Value  = "https://w3id.org/dpv/v2.0/use-cases#"
$ws.Cells.Item(7, 2).Value  = "https://w3id.org/dpv/v2.0/nace#"
$ws.Cells.Item(8, 2).Value  = "https://w3id.org/dpv/v2.0/legal/eu/gdpr#"
$ws.Cells.Item(9, 2).Value  = "https://w3id.org/dpv/v2.0/pd#"
$ws.Cells.Item(10, 2).Value = "https://w3id.org/dpv/v2.0/tech#"
$ws.Cells.Item(11, 2).Value = "https://w3id.org/dpv/v2.0/legal#"
$ws.Cells.Item(12, 2).Value = "https://w3id.org/dpv/v2.0/risk#"
$ws.Cells.Item(13, 2).Value = "https://w3id.org/dpv/v2.0/rights/eu#"
$ws.Cells.Item(14, 2).Value = "https://w3id.org/dpv/v2.0/legal/eu/dga#"
$ws.Cells.Item(15, 2).Value = "https://w3id.org/dpv/v2.0/legal/eu/aiact#"
$ws.Cells.Item(16, 2).Value = "https://w3id.org/dpv/v2.0/legal/eu/nis2#"
$ws.Cells.Item(17, 2).Value = "https://w3id.org/dpv/v2.0/loc#"
$ws.Cells.Item(18, 2).Value = "https://w3id.org/dpv/v2.0/legal/eu#"
$ws.Cells.Item(19, 2).Value = "https://w3id.org/dpv/v2.0/legal/de#"
$ws.Cells.Item(20, 2).Value = "https://w3id.org/dpv/v2.0/legal/ie#"
$ws.Cells.Item(21, 2).Value = "https://w3id.org/dpv/v2.0/legal/gb#"
$ws.Cells.Item(22, 2).Value = "https://w3id.org/dpv/v2.0/legal/us#"
$ws.Cells.Item(23, 2).Value = "https://w3id.org/dpv/v2.0/legal/in#"
$ws.Cells.Item(24, 2).Value = "https://w3id.org/dpv/v2.0/justifications#"
$ws.Cells.Item(25, 2).Value = "https://w3id.org/dpv/v2.0/ai#"

# Row 2's hyperlink style also moves from s=7 to s=10 (matching the rest
# of the column's link styling) per the diff.
$ws.Cells.Item(2, 2).Font.Underline = $true
$ws.Cells.Item(2, 2).Font.Color = $ws.Cells.Item(3, 2).Font.Color
